# Weekly fruit/vegetable price update:
# Insert two new rows (a new week's "Primera"/"Segunda" price records) right
# after the existing row 902, pushing the previously-recorded weeks down by
# two rows (old 903..941 -> new 905..943), and populate the two newly
# inserted rows with the latest week's data.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert 2 blank rows at position 903 (shifts rows 903:941 down to 905:943).
$ws.Rows.Item(903).Resize(2).Insert()

# ---- New row 903 ("Primera") ----
$ws.Cells.Item(903, 1).Value = 3
$ws.Cells.Item(903, 2).Value = "Femacal de La Calera"
$ws.Cells.Item(903, 3).Value = "Coquimbo"
$ws.Cells.Item(903, 4).Value = 44939
$ws.Cells.Item(903, 5).Value = 5
$ws.Cells.Item(903, 6).Value = 100114014
$ws.Cells.Item(903, 7).Value = "Betarraga"
$ws.Cells.Item(903, 8).Value = "Sin especificar"
$ws.Cells.Item(903, 9).Value = "Primera"
$ws.Cells.Item(903, 10).Value = 1800
$ws.Cells.Item(903, 11).Value = 500
$ws.Cells.Item(903, 12).Value = 500
$ws.Cells.Item(903, 13).Value = 500
$ws.Cells.Item(903, 14).Value = "`$/paquete 4 unidades"
$ws.Cells.Item(903, 15).Value = "Provincia de Quillota"
$ws.Cells.Item(903, 16).Value = 125
$ws.Cells.Item(903, 17).Value = 4
$ws.Cells.Item(903, 18).Value = "Hortaliza"

# ---- New row 904 ("Segunda") ----
$ws.Cells.Item(904, 1).Value = 3
$ws.Cells.Item(904, 2).Value = "Femacal de La Calera"
$ws.Cells.Item(904, 3).Value = "Coquimbo"
$ws.Cells.Item(904, 4).Value = 44939
$ws.Cells.Item(904, 5).Value = 5
$ws.Cells.Item(904, 6).Value = 100114014
$ws.Cells.Item(904, 7).Value = "Betarraga"
$ws.Cells.Item(904, 8).Value = "Sin especificar"
$ws.Cells.Item(904, 9).Value = "Segunda"
$ws.Cells.Item(904, 10).Value = 3400
$ws.Cells.Item(904, 11).Value = 350
$ws.Cells.Item(904, 12).Value = 400
$ws.Cells.Item(904, 13).Value = 376
$ws.Cells.Item(904, 14).Value = "`$/paquete 4 unidades"
$ws.Cells.Item(904, 15).Value = "Provincia de Quillota"
$ws.Cells.Item(904, 16).Value = 94
$ws.Cells.Item(904, 17).Value = 4
$ws.Cells.Item(904, 18).Value = "Hortaliza"
